$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.739.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.321.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.36%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.316.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.575"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "656.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.863.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.645.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.321.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.885"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("E23").Value = "  +5.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("E27").Value = "  -4.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "567.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.93%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.663.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.13%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0660"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.46%  "
